$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -11.7843
$ws.Range("B7").Value = 5.600199999999997
$ws.Range("D7").Value = -8.08279999999999
$ws.Range("A9").Value = -21.80870000000002
$ws.Range("D10").Value = -8.287299999999995
$ws.Range("B12").Value = 5.397899999999996
$ws.Range("D13").Value = -8.334400000000002
$ws.Range("B14").Value = 5.792100000000001
$ws.Range("C15").Value = -13.94369999999999
$ws.Range("D16").Value = -8.708100000000004
$ws.Range("A18").Value = -22.24550000000001
$ws.Range("A20").Value = -20.06269999999998
$ws.Range("D20").Value = -7.048599999999997
$ws.Range("D24").Value = -7.423500000000001
$ws.Range("B26").Value = 4.097700000000004
$ws.Range("A27").Value = -21.73139999999998
$ws.Range("B27").Value = 4.965300000000004
$ws.Range("B29").Value = 4.836499999999997
$ws.Range("C33").Value = -11.6066
$ws.Range("A35").Value = -19.9122
$ws.Range("C35").Value = -12.181
$ws.Range("B37").Value = 9.775500000000001
$ws.Range("B38").Value = 5.030600000000002
$ws.Range("C38").Value = -12.6972
$ws.Range("D39").Value = -7.316700000000005
$ws.Range("C43").Value = -13.5639
$ws.Range("C44").Value = -13.26579999999999
$ws.Range("C47").Value = -11.516
$ws.Range("D47").Value = -7.163399999999997
$ws.Range("D48").Value = -7.096399999999996
$ws.Range("B51").Value = 6.387500000000005
$ws.Range("C51").Value = -11.6178
$ws.Range("B52").Value = 5.1572
$ws.Range("D52").Value = -7.511299999999998
$ws.Range("B55").Value = 5.186899999999998
$ws.Range("D56").Value = -7.854799999999998
$ws.Range("C57").Value = -14.29799999999999
$ws.Range("C63").Value = -11.7257
$ws.Range("A69").Value = -21.68599999999999
$ws.Range("B69").Value = 5.280999999999994
$ws.Range("B70").Value = 6.298800000000004
$ws.Range("C70").Value = -11.6401
$ws.Range("A76").Value = -19.85429999999997
$ws.Range("A78").Value = -20.13339999999998
$ws.Range("B81").Value = 5.300200000000002
$ws.Range("A82").Value = -21.9387
$ws.Range("A83").Value = -21.88609999999999
$ws.Range("B83").Value = 7.119100000000003
$ws.Range("D84").Value = -8.979999999999997
$ws.Range("C88").Value = -11.87529999999999
$ws.Range("A93").Value = -20.62319999999998
$ws.Range("C99").Value = -12.4865
$ws.Range("D100").Value = -8.195099999999998
$ws.Range("D101").Value = -8.072699999999994
$ws.Range("B102").Value = 8.308400000000011
